$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data.
# A leading apostrophe forces Excel to store the value as literal text
# (matching the original inlineStr cells) instead of auto-converting
# number-looking strings (e.g. "116.60") into numeric values, which would
# otherwise strip formatting such as trailing zeros.
$ws.Range("D2").Value = "'52.282.94"
$ws.Range("E2").Value = "'  +1.98%  "
$ws.Range("D3").Value = "'2.797.59"
$ws.Range("E3").Value = "'  +1.79%  "
$ws.Range("E4").Value = "'  +0.03%  "
$ws.Range("D5").Value = "'347.15"
$ws.Range("E5").Value = "'  +4.34%  "
$ws.Range("D6").Value = "'116.60"
$ws.Range("E6").Value = "'  +1.20%  "
$ws.Range("D7").Value = "'0.552"
$ws.Range("E7").Value = "'  +4.14%  "
$ws.Range("E8").Value = "'  -0.05%  "
$ws.Range("D9").Value = "'0.594"
$ws.Range("E9").Value = "'  +3.78%  "
$ws.Range("D10").Value = "'43.04"
$ws.Range("E10").Value = "'  +3.98%  "
$ws.Range("D11").Value = "'0.0858"
$ws.Range("E11").Value = "'  +3.66%  "
$ws.Range("D12").Value = "'20.14"
$ws.Range("E12").Value = "'  -0.12%  "
$ws.Range("E13").Value = "'  +1.76%  "
$ws.Range("D14").Value = "'7.87"
$ws.Range("E14").Value = "'  +2.80%  "
$ws.Range("D15").Value = "'3.236.89"
$ws.Range("E15").Value = "'  +1.89%  "
$ws.Range("D16").Value = "'2.798.35"
$ws.Range("E16").Value = "'  +2.99%  "
$ws.Range("D17").Value = "'0.895"
$ws.Range("E17").Value = "'  +0.96%  "
$ws.Range("D18").Value = "'52.173.49"
$ws.Range("E18").Value = "'  +1.99%  "
$ws.Range("D19").Value = "'3.25"
$ws.Range("E19").Value = "'  +7.80%  "
$ws.Range("D20").Value = "'7.16"
$ws.Range("E20").Value = "'  +4.30%  "
$ws.Range("D21").Value = "'13.44"
$ws.Range("E21").Value = "'  -2.38%  "
$ws.Range("D22").Value = "'0.0₃0982"
$ws.Range("E22").Value = "'  +2.18%  "
$ws.Range("D23").Value = "'70.18"
$ws.Range("E23").Value = "'  +0.12%  "
$ws.Range("D24").Value = "'270.24"
$ws.Range("E24").Value = "'  -2.99%  "
$ws.Range("D25").Value = "'2.77"
$ws.Range("E25").Value = "'  +5.45%  "
$ws.Range("E26").Value = "'  -0.90%  "
$ws.Range("E27").Value = "'  -0.12%  "
$ws.Range("E28").Value = "'  -1.11%  "
$ws.Range("E29").Value = "'  +0.78%  "
$ws.Range("E30").Value = "'  -0.12%  "
$ws.Range("D31").Value = "'35.11"
$ws.Range("E31").Value = "'  -1.81%  "
$ws.Range("D32").Value = "'50.22"
$ws.Range("E32").Value = "'  +0.45%  "
$ws.Range("D33").Value = "'5.72"
$ws.Range("E33").Value = "'  +1.61%  "
$ws.Range("D34").Value = "'0.0435"
$ws.Range("E34").Value = "'  +23.84%  "
$ws.Range("D35").Value = "'0.0828"
$ws.Range("E35").Value = "'  +0.15%  "
$ws.Range("E36").Value = "'  +0.66%  "
$ws.Range("E37").Value = "'  +0.04%  "
$ws.Range("B38").Value = "'RenderToken"
$ws.Range("C38").Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").Value = "'5.00"
$ws.Range("E38").Value = "'  -1.07%  "
$ws.Range("B39").Value = "'Celestia"
$ws.Range("C39").Value = "'https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D39").Value = "'18.87"
$ws.Range("E39").Value = "'  -3.01%  "
$ws.Range("E40").Value = "'  +0.10%  "
$ws.Range("D41").Value = "'2.69"
$ws.Range("E41").Value = "'  +20.46%  "
$ws.Range("B42").Value = "'Monero"
$ws.Range("C42").Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D42").Value = "'127.73"
$ws.Range("E42").Value = "'  -1.31%  "
$ws.Range("B43").Value = "'EnergySwap"
$ws.Range("C43").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").Value = "'23.41"
$ws.Range("E43").Value = "'  -1.40%  "
$ws.Range("E44").Value = "'  +1.87%  "
$ws.Range("E45").Value = "'  +0.92%  "
$ws.Range("E46").Value = "'  -1.34%  "
$ws.Range("D47").Value = "'2.071.52"
$ws.Range("E47").Value = "'  -2.01%  "
$ws.Range("D48").Value = "'2.37"
$ws.Range("E48").Value = "'  +2.82%  "
$ws.Range("D49").Value = "'0.969"
$ws.Range("E49").Value = "'  +12.91%  "
$ws.Range("D50").Value = "'5.53"
$ws.Range("E50").Value = "'  -0.67%  "
$ws.Range("D51").Value = "'8.97"
$ws.Range("E51").Value = "'  -0.97%  "
